$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 in package order) - rows 3-13
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1017
$ws1.Range("F4").Value = 164
$ws1.Range("F5").Value = 2752
$ws1.Range("F6").Value = 92
$ws1.Range("F8").Value = 19
$ws1.Range("F9").Value = 116
$ws1.Range("F10").Value = 57
$ws1.Range("F11").Value = 49
$ws1.Range("F12").Value = 2561
$ws1.Range("F13").Value = 708

# Sheet "全部类型" (sheet4 in package order) - rows 4-15
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1017
$ws4.Range("F5").Value = 164
$ws4.Range("F6").Value = 2752
$ws4.Range("F7").Value = 92
$ws4.Range("F9").Value = 19
$ws4.Range("F11").Value = 116
$ws4.Range("F12").Value = 57
$ws4.Range("F13").Value = 49
$ws4.Range("F14").Value = 2561
$ws4.Range("F15").Value = 708
